$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.797.63"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "3.456.02"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Formula = "'582.37"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").Formula = "'146.94"
$ws.Range("E6").Value = "  +6.72%  "
$ws.Range("D7").Value = "3.456.72"
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D11").Formula = "'0.127"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("E12").Value = "  +2.51%  "
$ws.Range("D13").Value = "4.045.91"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Formula = "'27.90"
$ws.Range("E14").Value = "  +8.77%  "
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "3.456.17"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "61.905.25"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Formula = "'6.26"
$ws.Range("E19").Value = "  +8.75%  "
$ws.Range("D20").Formula = "'14.40"
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").Formula = "'9.56"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").Formula = "'389.39"
$ws.Range("E22").Value = "  +2.75%  "
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").Formula = "'73.62"
$ws.Range("E24").Value = "  +3.46%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Formula = "'0.999"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "3.598.30"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").Formula = "'0.180"
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Formula = "'7.73"
$ws.Range("E30").Value = "  +4.33%  "
$ws.Range("D31").Formula = "'0.998"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  -9.98%  "
$ws.Range("D33").Formula = "'8.23"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("E34").Value = "  +1.98%  "
$ws.Range("D36").Formula = "'24.18"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "3.485.52"
$ws.Range("E37").Value = "  +2.25%  "
$ws.Range("E38").Value = "  +2.29%  "
$ws.Range("E39").Value = "  +1.66%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Formula = "'166.76"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  +3.57%  "
$ws.Range("D43").Formula = "'27.29"
$ws.Range("E43").Value = "  +7.08%  "
$ws.Range("E44").Value = "  +4.16%  "
$ws.Range("D45").Formula = "'42.56"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("E46").Value = "  +3.45%  "
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("E49").Value = "  -1.91%  "
$ws.Range("D50").Value = "2.573.02"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("E51").Value = "  +1.96%  "
